# The deck ships two themes:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock Office palette), used by the Notes Master
#   ppt/theme/theme2.xml -> "Integral", used by the Slide Master (and therefore every slide)
#
# The edit swaps the two palettes: the Slide Master's theme becomes the stock
# "Office Theme" colours, while the "Integral" palette moves to the Notes
# Master's theme. In the PowerPoint object model the Slide Master's theme
# colours are reached through Master.Theme.ThemeColorScheme (a 12-entry
# RGBColor collection in the standard clrScheme order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), so we rewrite each entry in place with the
# "Office Theme" RGB values that used to live in theme1.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" colours (previously theme1.xml),
# written in clrScheme order. PowerPoint's RGB() builtin encodes a colour as
# R + G*256 + B*65536, so the literals below are that encoding of the
# corresponding "RRGGBB" hex values noted in each comment.
$officeThemeRGB = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}

# Best-effort: also try to rename the theme itself from "Integral" to
# "Office Theme" so the saved a:theme/a:clrScheme name attributes line up
# with the new palette. Wrapped defensively (and scoped to Theme.Name only -
# Design.Name is deliberately not used here because it renames the slide
# master's <p:cSld> instead of the theme, which is not part of this edit).
try { $master.Theme.Name = "Office Theme" } catch {}
